$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EVANGELINA's Saldo (row 2, column C): 1000000 -> 869000
$ws.Cells.Item(2, 3).Value = 869000

# Update DIEGO's Saldo (row 7, column C): 14082.35 -> 14000
$ws.Cells.Item(7, 3).Value = 14000

# Delete the MARCELO row (row 8, account 000772433 / 10000)
$ws.Rows.Item(8).Delete()

# After the MARCELO row is removed, the BRUNO row (was row 10, account
# 004584982 / 7012.11) shifts up to row 9 -- THOMAS (row 9) stayed put and is
# now row 8. The next three rows to delete (BRUNO, FLAVIA, ERICA) are now
# consecutive at row 9.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(9).Delete()
